$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

Set-TextValue "D2" "43.414.91"
Set-TextValue "E2" "  +0.92%  "
Set-TextValue "D3" "2.373.63"
Set-TextValue "E3" "  +3.00%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "310.57"
Set-TextValue "E5" "  +0.39%  "
Set-TextValue "D6" "105.10"
Set-TextValue "E6" "  +4.70%  "
Set-TextValue "D7" "0.523"
Set-TextValue "E7" "  -2.08%  "
Set-TextValue "E8" "  +0.07%  "
Set-TextValue "D9" "0.519"
Set-TextValue "E9" "  +1.88%  "
Set-TextValue "D10" "36.31"
Set-TextValue "E10" "  +0.97%  "
Set-TextValue "D11" "53.08"
Set-TextValue "E11" "  +2.01%  "
Set-TextValue "D12" "0.0814"
Set-TextValue "E12" "  -0.39%  "
Set-TextValue "D13" "0.113"
Set-TextValue "E13" "  -0.43%  "
Set-TextValue "D14" "7.01"
Set-TextValue "E14" "  +0.35%  "
Set-TextValue "D15" "2.742.00"
Set-TextValue "E15" "  +3.10%  "
Set-TextValue "D16" "15.64"
Set-TextValue "E16" "  +5.32%  "
Set-TextValue "D17" "2.371.33"
Set-TextValue "E17" "  +3.08%  "
Set-TextValue "D18" "0.816"
Set-TextValue "E18" "  +1.71%  "
Set-TextValue "D19" "43.377.84"
Set-TextValue "E19" "  +0.98%  "
Set-TextValue "D20" "12.04"
Set-TextValue "E20" "  -3.79%  "
Set-TextValue "E21" "  +0.49%  "
Set-TextValue "D22" "6.28"
Set-TextValue "E22" "  +3.28%  "
Set-TextValue "D23" "68.50"
Set-TextValue "E23" "  +0.49%  "
Set-TextValue "D24" "242.00"
Set-TextValue "E24" "  +0.85%  "
Set-TextValue "E25" "  +2.01%  "
Set-TextValue "D26" "2.63"
Set-TextValue "E26" "  +0.63%  "
Set-TextValue "E27" "  -0.23%  "
Set-TextValue "D28" "25.97"
Set-TextValue "E28" "  +7.25%  "
Set-TextValue "B29" "InjectiveProtocol"
Set-TextValue "C29" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D29" "36.91"
Set-TextValue "E29" "  -5.50%  "
Set-TextValue "B30" "Toncoin"
Set-TextValue "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "2.22"
Set-TextValue "E30" "  +4.57%  "
Set-TextValue "D31" "9.63"
Set-TextValue "E31" "  -0.01%  "
Set-TextValue "D32" "161.94"
Set-TextValue "E32" "  -4.49%  "
Set-TextValue "D33" "5.30"
Set-TextValue "E33" "  -0.57%  "
Set-TextValue "D34" "0.999"
Set-TextValue "E34" "  -0.07%  "
Set-TextValue "D35" "18.34"
Set-TextValue "E35" "  +3.44%  "
Set-TextValue "E36" "  +6.56%  "
Set-TextValue "D37" "3.14"
Set-TextValue "E37" "  -0.31%  "
Set-TextValue "D38" "4.70"
Set-TextValue "E38" "  +11.83%  "
Set-TextValue "B39" "Hedera"
Set-TextValue "C39" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D39" "0.0742"
Set-TextValue "E39" "  +0.34%  "
Set-TextValue "B40" "ARBITRUM"
Set-TextValue "C40" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "1.95"
Set-TextValue "E40" "  +5.83%  "
Set-TextValue "E41" "  +0.55%  "
Set-TextValue "D42" "0.114"
Set-TextValue "E42" "  -1.26%  "
Set-TextValue "D43" "2.47"
Set-TextValue "E43" "  +7.19%  "
Set-TextValue "D44" "20.32"
Set-TextValue "E44" "  +4.94%  "
Set-TextValue "D45" "2.002.21"
Set-TextValue "E45" "  +1.62%  "
Set-TextValue "E46" "  +0.41%  "
Set-TextValue "D47" "3.19"
Set-TextValue "E47" "  +6.24%  "
Set-TextValue "D48" "10.39"
Set-TextValue "E48" "  +6.21%  "
Set-TextValue "D49" "58.29"
Set-TextValue "E49" "  +5.85%  "
Set-TextValue "D50" "2.95"
Set-TextValue "E50" "  -1.49%  "
Set-TextValue "D51" "2.577.11"
Set-TextValue "E51" "  +1.97%  "
